$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 57.6
$ws.Range("I6").Value = 47.25
$ws.Range("K6").Value = 141.75
$ws.Range("M6").Value = -29.75
$ws.Range("H40").Value = 4255.2593
$ws.Range("J40").Value = 2991.4167
$ws.Range("L40").Value = 2991.4167
$ws.Range("N40").Value = -3341.4167
$ws.Range("H51").Value = 41303.4
$ws.Range("I51").Value = 36374.668
$ws.Range("K51").Value = 36374.668
$ws.Range("M51").Value = -35890.668
$ws.Range("H80").Value = 543.4
$ws.Range("I80").Value = 888.5
$ws.Range("J80").Value = 313.33334
$ws.Range("K80").Value = 2665.5
$ws.Range("L80").Value = 940.0000200000001
$ws.Range("M80").Value = -1667.5
$ws.Range("N80").Value = -2936.00002
$ws.Range("H83").Value = 543.4
$ws.Range("I83").Value = 888.5
$ws.Range("J83").Value = 313.33334
$ws.Range("K83").Value = 7996.5
$ws.Range("L83").Value = 2820.00006
$ws.Range("M83").Value = -3004.5
$ws.Range("N83").Value = -12804.00006
$ws.Range("H87").Value = 86967.164
$ws.Range("J87").Value = 94360.60000000001
$ws.Range("L87").Value = 94360.60000000001
$ws.Range("N87").Value = -96856.60000000001
$ws.Range("H90").Value = 86967.164
$ws.Range("J90").Value = 94360.60000000001
$ws.Range("L90").Value = 283081.8
$ws.Range("N90").Value = -295561.8
$ws.Range("H96").Value = 499.1
$ws.Range("I96").Value = 454.55554
$ws.Range("J96").Value = 900
$ws.Range("K96").Value = 1363.66662
$ws.Range("L96").Value = 2700
$ws.Range("M96").Value = 9.333380000000034
$ws.Range("N96").Value = -5446
$ws.Range("H98").Value = 1697.2858
$ws.Range("I98").Value = 1806.0968
$ws.Range("J98").Value = 854
$ws.Range("K98").Value = 1806.0968
$ws.Range("L98").Value = 854
$ws.Range("M98").Value = -308.0968
$ws.Range("N98").Value = -3850
$ws.Range("H115").Value = 3535
$ws.Range("I115").Value = 3535
$ws.Range("J115").Value = 0
$ws.Range("K115").Value = 10605
$ws.Range("L115").Value = 0
$ws.Range("M115").Value = -9038
$ws.Range("N115").ClearContents()
$ws.Range("H116").Value = 87187.5
$ws.Range("I116").Value = 157000
$ws.Range("J116").Value = 17375
$ws.Range("K116").Value = 157000
$ws.Range("L116").Value = 17375
$ws.Range("M116").Value = -153558
$ws.Range("N116").Value = -24259
$ws.Range("H122").Value = 1697.2858
$ws.Range("I122").Value = 1806.0968
$ws.Range("J122").Value = 854
$ws.Range("K122").Value = 5418.2904
$ws.Range("L122").Value = 2562
$ws.Range("M122").Value = -2968.2904
$ws.Range("N122").Value = -7462

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8780771
$ws.Range("I32").Value = 10210671
$ws.Range("K32").Value = 10210671
$ws.Range("M32").Value = -10210384
$ws.Range("H45").Value = 5050.9
$ws.Range("I45").Value = 4723.222
$ws.Range("K45").Value = 4723.222
$ws.Range("M45").Value = -4346.222
$ws.Range("H61").Value = 5947.8647
$ws.Range("I61").Value = 2502.7144
$ws.Range("K61").Value = 2502.7144
$ws.Range("M61").Value = -2290.7144
$ws.Range("H107").Value = 70000
$ws.Range("J107").Value = 70000
$ws.Range("L107").Value = 70000
$ws.Range("N107").Value = -77680
$ws.Range("H109").Value = 85000
$ws.Range("J109").Value = 85000
$ws.Range("L109").Value = 85000
$ws.Range("N109").Value = -87774
$ws.Range("H112").Value = 74176.55
$ws.Range("I112").Value = 83344.375
$ws.Range("J112").Value = 49729
$ws.Range("K112").Value = 83344.375
$ws.Range("L112").Value = 49729
$ws.Range("M112").Value = -81867.375
$ws.Range("N112").Value = -52683
$ws.Range("H132").Value = 780755.5600000001
$ws.Range("I132").Value = 1167060.2
$ws.Range("J132").Value = 8146.2
$ws.Range("K132").Value = 3501180.6
$ws.Range("L132").Value = 24438.6
$ws.Range("M132").Value = -3498650.6
$ws.Range("N132").Value = -29498.6
$ws.Range("H136").Value = 5947.8647
$ws.Range("I136").Value = 2502.7144
$ws.Range("K136").Value = 7508.1432
$ws.Range("M136").Value = -4958.1432

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 1987.25
$ws.Range("I107").Value = 1391.1818
$ws.Range("J107").Value = 3298.6
$ws.Range("K107").Value = 1391.1818
$ws.Range("L107").Value = 3298.6
$ws.Range("M107").Value = 528.8181999999999
$ws.Range("N107").Value = -7138.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H10").Value = 1141.6
$ws.Range("I10").Value = 702
$ws.Range("K10").Value = 2106
$ws.Range("M10").Value = -1967
$ws.Range("H60").Value = 1311.9565
$ws.Range("I60").Value = 1210.3846
$ws.Range("K60").Value = 3631.1538
$ws.Range("M60").Value = -3380.1538
$ws.Range("H61").Value = 149.33333
$ws.Range("I61").Value = 174
$ws.Range("J61").Value = 100
$ws.Range("K61").Value = 522
$ws.Range("L61").Value = 300
$ws.Range("M61").Value = -307
$ws.Range("N61").Value = -730
$ws.Range("H109").Value = 3582.5
$ws.Range("I109").Value = 1781.1666
$ws.Range("K109").Value = 5343.4998
$ws.Range("M109").Value = -4303.4998
$ws.Range("H129").Value = 2184.55
$ws.Range("I129").Value = 586.3333
$ws.Range("J129").Value = 2869.5
$ws.Range("K129").Value = 1758.9999
$ws.Range("L129").Value = 8608.5
$ws.Range("M129").Value = 3241.0001
$ws.Range("N129").Value = -18608.5
$ws.Range("H131").Value = 8877.843999999999
$ws.Range("J131").Value = 8877.843999999999
$ws.Range("L131").Value = 26633.532
$ws.Range("N131").Value = -36713.532
$ws.Range("H132").Value = 3888.9167
$ws.Range("I132").Value = 2998
$ws.Range("J132").Value = 4185.8887
$ws.Range("K132").Value = 26982
$ws.Range("L132").Value = 37672.99830000001
$ws.Range("M132").Value = -24452
$ws.Range("N132").Value = -42732.99830000001
$ws.Range("H137").Value = 2612.1333
$ws.Range("J137").Value = 2688
$ws.Range("L137").Value = 8064
$ws.Range("N137").Value = -18264
$ws.Range("H138").Value = 2064.1428
$ws.Range("I138").Value = 1289.8
$ws.Range("K138").Value = 3869.4
$ws.Range("M138").Value = 1270.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3433.8572
$ws.Range("I102").Value = 3082.6155
$ws.Range("K102").Value = 3082.6155
$ws.Range("M102").Value = -1460.6155
$ws.Range("H132").Value = 2445.3125
$ws.Range("I132").Value = 1738.6923
$ws.Range("K132").Value = 5216.0769
$ws.Range("M132").Value = -2686.0769

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3199.4
$ws.Range("I40").Value = 2856.2856
$ws.Range("J40").Value = 4000
$ws.Range("K40").Value = 2856.2856
$ws.Range("L40").Value = 4000
$ws.Range("M40").Value = -2720.2856
$ws.Range("N40").Value = -4272
$ws.Range("H68").Value = 9998.333000000001
$ws.Range("I68").Value = 9398
$ws.Range("K68").Value = 9398
$ws.Range("M68").Value = -8649
$ws.Range("H71").Value = 9998.333000000001
$ws.Range("I71").Value = 9398
$ws.Range("K71").Value = 46990
$ws.Range("M71").Value = -43246
$ws.Range("H82").Value = 4938.5386
$ws.Range("I82").Value = 3518.1667
$ws.Range("K82").Value = 3518.1667
$ws.Range("M82").Value = -3157.1667
$ws.Range("H85").Value = 4938.5386
$ws.Range("I85").Value = 3518.1667
$ws.Range("K85").Value = 3518.1667
$ws.Range("M85").Value = -2270.1667
$ws.Range("H110").Value = 77786.39999999999
$ws.Range("J110").Value = 77786.39999999999
$ws.Range("L110").Value = 77786.39999999999
$ws.Range("N110").Value = -85966.39999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 2957.9092
$ws.Range("I100").Value = 3087.6667
$ws.Range("K100").Value = 6175.3334
$ws.Range("M100").Value = -5634.3334
$ws.Range("H136").Value = 18140534
$ws.Range("I136").Value = 29301796
$ws.Range("J136").Value = 3484.875
$ws.Range("K136").Value = 87905388
$ws.Range("L136").Value = 10454.625
$ws.Range("M136").Value = -87902838
$ws.Range("N136").Value = -15554.625
